# ToDoList.xlsx update: rename a task, add a new task row, and fill in
# the "Result" column notes for the GIT/임베디드 section task block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToDo")

# Rename the "Virtual Box 환경설정 ST" task (B16) to reflect the updated scope.
$ws.Range("B16").Value = "Virtual Box 설치,환경설정 ST"

# The ".so library ST" task that used to live at B18 moves down to B20 ...
$ws.Range("B20").Value = ".so library ST"

# ... and B18 becomes a brand new task about the Ubuntu VM setup.
$ws.Range("B18").Value = "Virtual Box Ubuntu 설치, 환경설정ST"

# Fill in the "Result" column (column C) notes for this block of tasks.
$ws.Range("C14").Value = "SS BBUe"
$ws.Range("C16").Value = "Init 구조 ST"
$ws.Range("C17").Value = "GOM FPGA, Device Driver"
$ws.Range("C18").Value = "Installation UI"

# Leave the last active selection where the author ended up.
$ws.Range("C19").Select()
